$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change header from 'barcode' to 'product_id' (single column now accepts
# either a SKU or a barcode value for the cascading product search).
$ws.Range("A1").Value = "product_id"

# Update column A values to the new mixed SKU/barcode identifiers.
# Values that are entirely digits are entered with a leading apostrophe
# so Excel keeps them as text (matching how the original barcodes, which
# also look numeric, were stored as text in the workbook).
$ws.Range("A2").Value = "SKU123456"
$ws.Range("A3").Value = "'4607034370244"
$ws.Range("A4").Value = "OZON789012"
$ws.Range("A5").Value = "'9999999999999"
$ws.Range("A6").Value = "UNKNOWN_SKU"
